$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.766.01'
$ws.Range("E2").Value = '  -4.54%  '
$ws.Range("D3").Value = '2.988.95'
$ws.Range("E3").Value = '  -4.64%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '542.27'
$ws.Range("E5").Value = '  -5.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.74'
$ws.Range("E6").Value = '  -7.22%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.570'
$ws.Range("E8").Value = '  -0.50%  '
$ws.Range("D9").Value = '3.001.51'
$ws.Range("E9").Value = '  -4.61%  '
$ws.Range("E10").Value = '  -3.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.14'
$ws.Range("E11").Value = '  -7.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.369'
$ws.Range("E12").Value = '  -3.08%  '
$ws.Range("D13").Value = '3.511.32'
$ws.Range("E13").Value = '  -4.71%  '
$ws.Range("E14").Value = '  -1.02%  '
$ws.Range("D15").Value = '61.828.73'
$ws.Range("E15").Value = '  -4.41%  '
$ws.Range("E16").Value = '  -3.69%  '
$ws.Range("D17").Value = '2.998.02'
$ws.Range("E17").Value = '  -4.53%  '
$ws.Range("E18").Value = '  -5.17%  '
$ws.Range("E19").Value = '  -1.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.09'
$ws.Range("E20").Value = '  -2.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '378.87'
$ws.Range("E21").Value = '  -8.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.71'
$ws.Range("E22").Value = '  -4.51%  '
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("E24").Value = '  -3.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.12'
$ws.Range("E25").Value = '  -3.62%  '
$ws.Range("D26").Value = '3.113.76'
$ws.Range("E26").Value = '  -4.79%  '
$ws.Range("E27").Value = '  -2.38%  '
$ws.Range("E28").Value = '  -3.10%  '
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("D30").Value = '0.0₃0936'
$ws.Range("E30").Value = '  -8.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.24'
$ws.Range("E31").Value = '  -8.98%  '
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("E33").Value = '  -4.58%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '20.45'
$ws.Range("E34").Value = '  -3.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '161.11'
$ws.Range("E35").Value = '  -1.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.98'
$ws.Range("E36").Value = '  -3.90%  '
$ws.Range("E37").Value = '  -4.76%  '
$ws.Range("E38").Value = '  -4.94%  '
$ws.Range("E39").Value = '  -5.35%  '
$ws.Range("E40").Value = '  -7.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '37.68'
$ws.Range("E41").Value = '  -1.51%  '
$ws.Range("D42").Value = '2.418.43'
$ws.Range("E42").Value = '  -7.60%  '
$ws.Range("E43").Value = '  -5.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.08'
$ws.Range("E44").Value = '  -6.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.672'
$ws.Range("E45").Value = '  -2.65%  '
$ws.Range("E46").Value = '  -3.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.18'
$ws.Range("E47").Value = '  -1.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.997'
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("E49").Value = '  -3.66%  '
$ws.Range("E50").Value = '  -2.25%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.76'
$ws.Range("E51").Value = '  -6.55%  '
